# feat: add 2022-Q1 data
#
# Inserts a new "2022-Q1" worksheet (between "2021-Q4" and "总计") with the
# quarterly fund-holdings detail, and updates the "总计" (totals) sheet with
# a new leading row summarising the 2022-Q1 quarter.

$wb = $excel.ActiveWorkbook

$wsQ4 = $wb.Worksheets.Item("2021-Q4")
$wsOldTotal = $wb.Worksheets.Item("总计")

# Drop the existing "总计" sheet - it gets rebuilt after "2022-Q1" so that
# the sheet order (and freed-up internal sheetId) comes out right: the new
# "2022-Q1" tab takes the id vacated by "总计", and "总计" is re-created
# after it.
$wsOldTotal.Delete()

# ---------------------------------------------------------------------
# New "2022-Q1" sheet, positioned right after "2021-Q4"
# ---------------------------------------------------------------------
$wsQ1 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsQ4)
$wsQ1.Name = "2022-Q1"

# Copy the header-row and "index column" formatting from the "2021-Q4"
# sheet (same bold/border/centred style used on every per-quarter sheet).
$wsQ4.Range("B1:H1").Copy()
$wsQ1.Range("B1:H1").PasteSpecial(-4122)
$wsQ4.Range("A2").Copy()
$wsQ1.Range("A2:A11").PasteSpecial(-4122)

$wsQ1.Range("B1").Value = "基金代码"
$wsQ1.Range("C1").Value = "基金名称"
$wsQ1.Range("D1").Value = "基金规模"
$wsQ1.Range("E1").Value = "股票总仓位"
$wsQ1.Range("F1").Value = "仓位占比"
$wsQ1.Range("G1").Value = "持有市值(亿元)"
$wsQ1.Range("H1").Value = "仓位排名"

# code, name, scale, position, ratio, value, rank
$q1Rows = @(
    @("001410", "信达澳银新能源产业股票",         "140.41", "92.06", "1.44", "2.0219", 8),
    @("012608", "信达澳银领先智选混合型证券投资基金", "38.78", "90.57", "1.39", "0.5390", 10),
    @("006257", "信达澳银先进智造股票",             "21.53", "93.84", "1.54", "0.3316", 8),
    @("009511", "信达澳银研究优选混合",             "9.41",  "92.12", "1.39", "0.1308", 10),
    @("009055", "圆信永丰大湾区主题混合A",          "1.60",  "92.94", "3.32", "0.0531", 7),
    @("009056", "圆信永丰大湾区主题混合C",          "1.24",  "92.94", "3.32", "0.0412", 7),
    @("310318", "申万菱信沪深300指数增强A",         "7.60",  "90.31", "0.05", "0.0038", 9),
    @("004976", "华润元大景泰混合A",               "1.79",  "37.61", "0.11", "0.0020", 9),
    @("004977", "华润元大景泰混合C",               "1.79",  "37.61", "0.11", "0.0020", 9),
    @("007804", "申万菱信沪深300指数增强C",         "0.87",  "90.31", "0.05", "0.0004", 9)
)

# Columns B, D, E, F, G hold numeric-looking text (fund codes with leading
# zeros, percentages, NAV figures) that must stay TEXT - exactly like the
# source workbook (generated with inline strings, not numbers). A leading
# apostrophe is how Excel enters a number-shaped value as text.
$rowNum = 2
foreach ($fund in $q1Rows) {
    $wsQ1.Cells.Item($rowNum, 1).Value = $rowNum - 2
    $wsQ1.Cells.Item($rowNum, 2).Value = "'" + $fund[0]
    $wsQ1.Cells.Item($rowNum, 3).Value = $fund[1]
    $wsQ1.Cells.Item($rowNum, 4).Value = "'" + $fund[2]
    $wsQ1.Cells.Item($rowNum, 5).Value = "'" + $fund[3]
    $wsQ1.Cells.Item($rowNum, 6).Value = "'" + $fund[4]
    $wsQ1.Cells.Item($rowNum, 7).Value = "'" + $fund[5]
    $wsQ1.Cells.Item($rowNum, 8).Value = $fund[6]
    $rowNum = $rowNum + 1
}

# ---------------------------------------------------------------------
# Re-create the "总计" sheet right after "2022-Q1", with a new leading
# row for the 2022-Q1 totals.
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsQ1)
$wsTotal.Name = "总计"

$wsQ4.Range("B1:D1").Copy()
$wsTotal.Range("B1:D1").PasteSpecial(-4122)
$wsQ4.Range("A2").Copy()
$wsTotal.Range("A2:A4").PasteSpecial(-4122)

$wsTotal.Range("B1").Value = "日期"
$wsTotal.Range("C1").Value = "持有数量(只)"
$wsTotal.Range("D1").Value = "持有市值(亿元)"

$allTotals = @(
    @("2022-Q1", 10, 3.13),
    @("2021-Q4", 3, 0.54),
    @("2021-Q3", 4, 0)
)

$rowNum = 2
foreach ($t in $allTotals) {
    $wsTotal.Cells.Item($rowNum, 1).Value = $rowNum - 2
    $wsTotal.Cells.Item($rowNum, 2).Value = $t[0]
    $wsTotal.Cells.Item($rowNum, 3).Value = $t[1]
    $wsTotal.Cells.Item($rowNum, 4).Value = $t[2]
    $rowNum = $rowNum + 1
}
